$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "65.579.53", "  -3.19%  ")
  3 = @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.793.81", "  +0.83%  ")
  4 = @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.12%  ")
  5 = @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "419.27", "  -0.46%  ")
  6 = @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "126.71", "  -4.32%  ")
  7 = @("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.791.85", "  +1.11%  ")
  8 = @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.597", "  -8.48%  ")
  9 = @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.999", "  -0.08%  ")
  10 = @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.708", "  -9.01%  ")
  11 = @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.161", "  -14.11%  ")
  12 = @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000340", "  -21.06%  ")
  13 = @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "39.63", "  -7.84%  ")
  14 = @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.426.28", "  +1.64%  ")
  15 = @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "9.79", "  -6.04%  ")
  16 = @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "15.63", "  +19.43%  ")
  17 = @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.137", "  -1.75%  ")
  18 = @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.800.80", "  +1.27%  ")
  19 = @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "19.22", "  -7.00%  ")
  20 = @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "66.064.59", "  -2.42%  ")
  21 = @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.05", "  -8.40%  ")
  22 = @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "396.90", "  -12.15%  ")
  23 = @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "14.08", "  -11.60%  ")
  24 = @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "82.97", "  -7.27%  ")
  25 = @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.95", "  -5.02%  ")
  26 = @("LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "5.74", "  +12.53%  ")
  27 = @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "36.49", "  -6.57%  ")
  28 = @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.16", "  -5.18%  ")
  29 = @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "9.22", "  -9.17%  ")
  30 = @("Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "687.96", "  +0.27%  ")
  31 = @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.72", "  -1.76%  ")
  32 = @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.119", "  -5.62%  ")
  33 = @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "12.10", "  -5.02%  ")
  34 = @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "7.22", "  -0.59%  ")
  35 = @("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.148", "  -10.61%  ")
  36 = @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.998", "  -0.15%  ")
  37 = @("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "37.29", "  -11.56%  ")
  38 = @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "54.43", "  -4.78%  ")
  39 = @("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0734", "  -4.33%  ")
  40 = @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0444", "  -10.26%  ")
  41 = @("ThetaToken", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", "2.88", "  -2.96%  ")
  42 = @("FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.01", "  +0.94%  ")
  43 = @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.133", "  -11.22%  ")
  44 = @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "4.39", "  +0.92%  ")
  45 = @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "143.38", "  -3.19%  ")
  46 = @("ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "3.07", "  -3.14%  ")
  47 = @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "3.21", "  -6.17%  ")
  48 = @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "2.02", "  -5.93%  ")
  49 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "25.68", "  -8.01%  ")
  50 = @("WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.52", "  -4.75%  ")
  51 = @("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "2.70", "  -7.48%  ")
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 4).NumberFormat = "@"
  $ws.Cells.Item($row, 4).Value = $vals[2]
  $ws.Cells.Item($row, 5).Value = $vals[3]
}
